# Daily update at 8 AM UTC
# Adds the next day's row (day 45620) to the "Wins Over Time" sheet and
# restores the previous last row (35) to the regular datetime number format,
# since row 36 is now the newest/last row and takes the "date only" format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 35 is no longer the last row, so it goes back to the standard
# "YYYY-MM-DD HH:MM:SS" format used by every other non-final row.
$ws.Cells.Item(35, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New last row (36): date-only format on column A, like every prior last row.
$ws.Cells.Item(36, 1).Value = 45620
$ws.Cells.Item(36, 1).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(36, 2).Value = 92
$ws.Cells.Item(36, 3).Value = 73
$ws.Cells.Item(36, 4).Value = 85
